# Generate Report for Handback
# - Flip "Ready for handoff" -> "Handed back: in sync with en-US" everywhere it appears
#   (Overview sheet zh-cn/de-de status columns, and Status column on the zh-cn / de-de sheets)
# - Fill in "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
#   for both rows on the zh-cn and de-de sheets, now that handback has happened
# - "Latest Target File" becomes a hyperlink back to the source .md file, just like
#   the existing "Source File Name" column link

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$handedBack = "Handed back: in sync with en-US"

# ----- Overview sheet: zh-cn / de-de status columns -----
$wsOverview.Range("E2").Value = $handedBack
$wsOverview.Range("F2").Value = $handedBack
$wsOverview.Range("E3").Value = $handedBack
$wsOverview.Range("F3").Value = $handedBack

# ----- zh-cn sheet -----
$wsZhCn.Range("C2").Value = $handedBack
$wsZhCn.Range("C3").Value = $handedBack

$wsZhCn.Range("I2").Value = "8aa2567e-b409-458a-a9ea-f8c40dd83391.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8e861ea329665b9b5b0879684cabd4ecd9939d1e/e2e/8aa2567e-b409-458a-a9ea-f8c40dd83391.md", "", "", "8aa2567e-b409-458a-a9ea-f8c40dd83391.md") | Out-Null
$wsZhCn.Range("J2").Value = "8aa2567e-b409-458a-a9ea-f8c40dd83391.30bd3b7c8efa0f089fd1219ee418dfebbd30e816.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-10-17 13:49:22"

$wsZhCn.Range("I3").Value = "f5348948-4f86-4d36-b8a0-67a8c1d3ffcf.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8e861ea329665b9b5b0879684cabd4ecd9939d1e/e2e/f5348948-4f86-4d36-b8a0-67a8c1d3ffcf.md", "", "", "f5348948-4f86-4d36-b8a0-67a8c1d3ffcf.md") | Out-Null
$wsZhCn.Range("J3").Value = "f5348948-4f86-4d36-b8a0-67a8c1d3ffcf.caacfc6d75049a83d9b85a9a59f4b3723427d218.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-10-17 13:49:22"

# ----- de-de sheet -----
$wsDeDe.Range("C2").Value = $handedBack
$wsDeDe.Range("C3").Value = $handedBack

$wsDeDe.Range("I2").Value = "8aa2567e-b409-458a-a9ea-f8c40dd83391.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8e861ea329665b9b5b0879684cabd4ecd9939d1e/e2e/8aa2567e-b409-458a-a9ea-f8c40dd83391.md", "", "", "8aa2567e-b409-458a-a9ea-f8c40dd83391.md") | Out-Null
$wsDeDe.Range("J2").Value = "8aa2567e-b409-458a-a9ea-f8c40dd83391.30bd3b7c8efa0f089fd1219ee418dfebbd30e816.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-10-17 13:49:59"

$wsDeDe.Range("I3").Value = "f5348948-4f86-4d36-b8a0-67a8c1d3ffcf.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8e861ea329665b9b5b0879684cabd4ecd9939d1e/e2e/f5348948-4f86-4d36-b8a0-67a8c1d3ffcf.md", "", "", "f5348948-4f86-4d36-b8a0-67a8c1d3ffcf.md") | Out-Null
$wsDeDe.Range("J3").Value = "f5348948-4f86-4d36-b8a0-67a8c1d3ffcf.caacfc6d75049a83d9b85a9a59f4b3723427d218.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-10-17 13:49:59"

Write-Output "Generated handback report"
